# Applies the edit described by the commit "My report version added. And
# new results included." to the single-slide poster deck.
#
# The underlying author edit was a light copy-edit pass over several body
# textboxes (re-typing/retouching text that PowerPoint had previously split
# into many same-formatted runs - the wording itself is unchanged, only the
# run boundaries collapse down to one run per contiguous span) plus a resize
# / reposition of the green "Conclusion" banner and its label.

function Merge-RunsInShape($shape, [string]$searchText) {
    # Re-assigns the exact same text back onto the character range that
    # currently spans $searchText. PowerPoint's text engine collapses every
    # run fully covered by the written range into a single run (taking on
    # the run-properties of the first run in the range), exactly mirroring
    # what happens when an author retypes/edits text that used to be split
    # across many runs with identical formatting.
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    $startIdx = $full.IndexOf($searchText)
    if ($startIdx -lt 0) {
        Write-Host "WARNING: text not found in shape $($shape.Name): $searchText"
        return
    }
    $len = $searchText.Length
    $sub = $tr.Characters($startIdx + 1, $len)
    $sub.Text = $searchText
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- TextBox 67 ("An investigation was undertaken ...") ---
$tb67 = $s.Shapes.Item(2)
Merge-RunsInShape $tb67 "The investigation narrowed down the scope of computer usage to web-browsing. The techniques seen to make web-browsing simpler can be extended to make computers generally more usable for the elderly."
Merge-RunsInShape $tb67 "Voice recognition was perceived as one of these techniques. Thus methods of voice referencing and visual annotations were specifically investigated."
Merge-RunsInShape $tb67 "Determine which of two voice referencing techniques perform better as a means of referencing links (accuracy). The referencing  techniques  used are:"

# --- TextBox 46 ("A facsimile of a news website ...") ---
$tb46 = $s.Shapes.Item(11)
Merge-RunsInShape $tb46 "A facsimile of a news website was designed. Numerical and link name referencing styles were applied to the website to investigate the performance of these techniques on a more complex web application. The preferred feedback techniques (highlighting and verbal feedback) were also incorporated into the site. Additional usage questions were also posed. "
Merge-RunsInShape $tb46 "  require  internet access to process commands. "

# --- TextBox 41 ("A simple website composed of questions and answers ...") ---
$tb41 = $s.Shapes.Item(12)
Merge-RunsInShape $tb41 " However the results in figure 2 indicates that link name "
Merge-RunsInShape $tb41 " Figure 3 indicates that users prefer link highlighting as a visual"

# --- TextBox 42 ("The large error rate observed ...") ---
$tb42 = $s.Shapes.Item(18)
Merge-RunsInShape $tb42 "The large error rate observed for numerical referencing (see Analysis) in Iteration 1 could largely be attributed to the user "
Merge-RunsInShape $tb42 " period. For this reason, the test was restructured. Feedback techniques were eliminated and a warm-up period (tutorial) for both numerical and spoken link name testing was provided."

# --- Resize/reposition the green "Conclusion" banner ---
$rect36 = $s.Shapes.Item(54)
$rect36.Width = 671.28585
$rect36.Height = 77.50524

$textbox37 = $s.Shapes.Item(55)
$textbox37.Left = 1618.88366
